# Realestate: append the 2023-05-29 10:27 resale-number snapshot as a new
# row (row 6) at the bottom of the CityResaleNum sheet, matching the
# existing layout (Date/Time/Weekday/Week stored as text, city figures as
# numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 6

# Style of an existing data row, used to keep the new row's formatting
# identical to its neighbours (no special number format).
$defaultStyle = $ws.Cells.Item(2, 1).Style

# Columns A and D hold values that look numeric/date-like ("2023-05-29",
# "22") but must stay plain text, just like the rest of the sheet. A
# leading apostrophe forces Excel to store them as text instead of
# auto-converting to a date serial / number; re-applying the plain data
# row's style afterwards strips the resulting quote-prefix formatting so
# the cell ends up styled exactly like its neighbours.
$ws.Cells.Item($row, 1).Value = "'2023-05-29"
$ws.Cells.Item($row, 1).Style = $defaultStyle

$ws.Cells.Item($row, 2).Value = "10:24:58"
$ws.Cells.Item($row, 3).Value = "Monday"

$ws.Cells.Item($row, 4).Value = "'22"
$ws.Cells.Item($row, 4).Style = $defaultStyle

$ws.Cells.Item($row, 5).Value = 119562
$ws.Cells.Item($row, 6).Value = 133603
$ws.Cells.Item($row, 7).Value = 157812
$ws.Cells.Item($row, 8).Value = 130683
$ws.Cells.Item($row, 9).Value = 174324
$ws.Cells.Item($row, 10).Value = 113280
$ws.Cells.Item($row, 11).Value = 198205
$ws.Cells.Item($row, 12).Value = 219955
$ws.Cells.Item($row, 13).Value = 172219
$ws.Cells.Item($row, 14).Value = 119843
$ws.Cells.Item($row, 15).Value = 38587
$ws.Cells.Item($row, 16).Value = 34923
$ws.Cells.Item($row, 17).Value = 50370
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36790
$ws.Cells.Item($row, 20).Value = -1

$wb.Save()
